$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Candidate credentials
$ws.Range("A2").Value = 'HqFKg679'
$ws.Range("B2").Value = 23110947
$ws.Range("C2").Value = 'sfjauje16'
$ws.Range("D2").Value = 'fr87#BV$'
$ws.Range("F2").Value = 'rnTynNPp'
$ws.Range("G2").Value = 'reWS'

# Row 3 - Candidate credentials
$ws.Range("A3").Value = 'dlkfu240'
$ws.Range("B3").Value = 23110946
$ws.Range("C3").Value = 'kjqfxwa63'
$ws.Range("D3").Value = 'QJ5f2&%u'
$ws.Range("F3").Value = 'pgxTBvqY'
$ws.Range("G3").Value = 'xNPm'
